$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells that already carry the two checkmark styles used throughout
# the sheet: a "no fill" checkmark style (like C3) and a "gray fill" checkmark
# style (like D3). We copy their formatting onto the cells that the diff
# turns into checkmarks, then write the checkmark character itself.

$plainCheck = $ws.Range("C3")   # Wingdings checkmark, no fill (style used by e.g. C3/E3)
$grayCheck  = $ws.Range("D3")   # Wingdings checkmark, gray fill (style used by e.g. D3/F3)

# F5: blank gray cell -> checkmark on gray fill
$grayCheck.Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = "ü"

# E12: blank cell -> checkmark, no fill
$plainCheck.Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = "ü"

# E20: blank cell -> checkmark, no fill
$plainCheck.Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "ü"

# G21: blank cell -> checkmark, no fill
$plainCheck.Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("G21").Value = "ü"

$excel.CutCopyMode = 0
